{"js": "const body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\nconst lastPara = paras.items[paras.items.length - 1];\n\n// Sequence of runs to append, mirroring the authored R code block:\n// a line break, then the new \"Baltimore <- Baltimore %>%\" pipe line,\n// then a mutate(season = case_when(...)) call spanning several lines.\nconst items = [\n  { br: true },\n  { text: \"Baltimore \", style: \"NormalTok\" },\n  { text: \"<-\", style: \"OtherTok\" },\n  { text: \" Baltimore \", style: \"NormalTok\" },\n  { text: \"%>%\", style: \"SpecialCharTok\" },\n  { text: \" \", style: \"NormalTok\" },\n  { br: true },\n  { text: \"  \", style: \"NormalTok\" },\n  { text: \"mutate\", style: \"FunctionTok\" },\n  { text: \"(\", style: \"NormalTok\" },\n  { text: \"season =\", style: \"AttributeTok\" },\n  { text: \" \", style: \"NormalTok\" },\n  { text: \"case_when\", style: \"FunctionTok\" },\n  { text: \"(month \", style: \"NormalTok\" },\n  { text: \">=\", style: \"SpecialCharTok\" },\n  { text: \"5\", style: \"DecValTok\" },\n  { text: \" \", style: \"NormalTok\" },\n  { text: \"&\", style: \"SpecialCharTok\" },\n  { text: \" month \", style: \"NormalTok\" },\n  { text: \"<=\", style: \"SpecialCharTok\" },\n  { text: \"10\", style: \"DecValTok\" },\n  { text: \" \", style: \"NormalTok\" },\n  { text: \"~\", style: \"SpecialCharTok\" },\n  { text: \"'Summer'\", style: \"StringTok\" },\n  { text: \",\", style: \"NormalTok\" },\n  { br: true },\n  { text: \"                            month \", style: \"NormalTok\" },\n  { text: \"<\", style: \"SpecialCharTok\" },\n  { text: \"5\", style: \"DecValTok\" },\n  { text: \" \", style: \"NormalTok\" },\n  { text: \"~\", style: \"SpecialCharTok\" },\n  { text: \" \", style: \"NormalTok\" },\n  { text: \"'Winter'\", style: \"StringTok\" },\n  { text: \",\", style: \"NormalTok\" },\n  { br: true },\n  { text: \"                            month \", style: \"NormalTok\" },\n  { text: \">\", style: \"SpecialCharTok\" },\n  { text: \"10\", style: \"DecValTok\" },\n  { text: \" \", style: \"NormalTok\" },\n  { text: \"~\", style: \"SpecialCharTok\" },\n  { text: \"'Winter'\", style: \"StringTok\" },\n  { text: \"))\", style: \"NormalTok\" },\n];\n\nfor (const item of items) {\n  const tail = lastPara.getRange(\"End\");\n  if (item.br) {\n    tail.insertBreak(Word.BreakType.line, Word.InsertLocation.end);\n  } else {\n    const inserted = tail.insertText(item.text, Word.InsertLocation.replace);\n    inserted.style = item.style;\n  }\n  await context.sync();\n}", "ps1": "$d = $word.ActiveDocument\n$p = $d.Paragraphs.Last\n\n# Helper: append $text to the end of $para, then (re-)apply the given\n# character style (rStyle) to just the text that was inserted. The style\n# range is computed from the paragraph length *before* the insertion,\n# minus 1, because this document's last paragraph has no trailing\n# paragraph mark (it ends exactly at the end of the document body), so\n# \"end of paragraph\" and \"end of document\" coincide.\nfunction Append-StyledText($para, $text, $styleName) {\n    $before = $para.Range.End\n    $para.Range.InsertAfter($text)\n    $start = $before - 1\n    $end = $start + $text.Length\n    $styled = $d.Range($start, $end)\n    $styled.Style = $styleName\n}\n\n# Manual line break == Chr(11) (vertical tab), which Word stores as <w:br/>.\n# Inserted directly (not passed through a function parameter) because the\n# control character does not survive a function-parameter hand-off here.\n\n$p.Range.InsertAfter([char]11)\nAppend-StyledText $p \"Baltimore \" \"NormalTok\"\nAppend-StyledText $p \"<-\" \"OtherTok\"\nAppend-StyledText $p \" Baltimore \" \"NormalTok\"\nAppend-StyledText $p \"%>%\" \"SpecialCharTok\"\nAppend-StyledText $p \" \" \"NormalTok\"\n$p.Range.InsertAfter([char]11)\nAppend-StyledText $p \"  \" \"NormalTok\"\nAppend-StyledText $p \"mutate\" \"FunctionTok\"\nAppend-StyledText $p \"(\" \"NormalTok\"\nAppend-StyledText $p \"season =\" \"AttributeTok\"\nAppend-StyledText $p \" \" \"NormalTok\"\nAppend-StyledText $p \"case_when\" \"FunctionTok\"\nAppend-StyledText $p \"(month \" \"NormalTok\"\nAppend-StyledText $p \">=\" \"SpecialCharTok\"\nAppend-StyledText $p \"5\" \"DecValTok\"\nAppend-StyledText $p \" \" \"NormalTok\"\nAppend-StyledText $p \"&\" \"SpecialCharTok\"\nAppend-StyledText $p \" month \" \"NormalTok\"\nAppend-StyledText $p \"<=\" \"SpecialCharTok\"\nAppend-StyledText $p \"10\" \"DecValTok\"\nAppend-StyledText $p \" \" \"NormalTok\"\nAppend-StyledText $p \"~\" \"SpecialCharTok\"\nAppend-StyledText $p \"'Summer'\" \"StringTok\"\nAppend-StyledText $p \",\" \"NormalTok\"\n$p.Range.InsertAfter([char]11)\nAppend-StyledText $p \"                            month \" \"NormalTok\"\nAppend-StyledText $p \"<\" \"SpecialCharTok\"\nAppend-StyledText $p \"5\" \"DecValTok\"\nAppend-StyledText $p \" \" \"NormalTok\"\nAppend-StyledText $p \"~\" \"SpecialCharTok\"\nAppend-StyledText $p \" \" \"NormalTok\"\nAppend-StyledText $p \"'Winter'\" \"StringTok\"\nAppend-StyledText $p \",\" \"NormalTok\"\n$p.Range.InsertAfter([char]11)\nAppend-StyledText $p \"                            month \" \"NormalTok\"\nAppend-StyledText $p \">\" \"SpecialCharTok\"\nAppend-StyledText $p \"10\" \"DecValTok\"\nAppend-StyledText $p \" \" \"NormalTok\"\nAppend-StyledText $p \"~\" \"SpecialCharTok\"\nAppend-StyledText $p \"'Winter'\" \"StringTok\"\nAppend-StyledText $p \"))\" \"NormalTok\"\n"}
